$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New FedEx tracking numbers (shared-string values) for the "ShipmentTrackNum"
# (column C) and "PackageTrackNum" (column D) columns of rows 2-22.
# Writing via .Value directly would make the engine infer a Number type for
# these purely-numeric strings (and forcing Text via NumberFormat/quote-prefix
# mints a brand new cell style, which the target file does not have).
# Instead we compute each value as TEXT() in an out-of-the-way scratch cell,
# copy it, and paste-special (values only) into the destination - this keeps
# the destination cell's existing (default, style-less) formatting untouched
# while still storing the result as a shared string, matching the diff.

$updates = @(
  @{ Cell = "C2";  Value = "320018624657" },
  @{ Cell = "C3";  Value = "320018621073" },
  @{ Cell = "C4";  Value = "320018621100" },
  @{ Cell = "C5";  Value = "320018621121" },
  @{ Cell = "C6";  Value = "320018621165" },
  @{ Cell = "C7";  Value = "320018621187" },
  @{ Cell = "C8";  Value = "320018621213" },
  @{ Cell = "C9";  Value = "320018621235" },
  @{ Cell = "C10"; Value = "320018621268" },
  @{ Cell = "C11"; Value = "320018621280" },
  @{ Cell = "C12"; Value = "320018621327" },
  @{ Cell = "C13"; Value = "320018621349" },
  @{ Cell = "C14"; Value = "320018621371" },
  @{ Cell = "C15"; Value = "320018621393" },
  @{ Cell = "C16"; Value = "320018621420" },
  @{ Cell = "C17"; Value = "320018621441" },
  @{ Cell = "C18"; Value = "320018621485" },
  @{ Cell = "C19"; Value = "320018621500" },
  @{ Cell = "C20"; Value = "320018621533" },
  @{ Cell = "C21"; Value = "320018621555" },
  @{ Cell = "C22"; Value = "320018621588" },
  @{ Cell = "D5";  Value = "320018621121" },
  @{ Cell = "D6";  Value = "320018621165" },
  @{ Cell = "D7";  Value = "320018621187" },
  @{ Cell = "D13"; Value = "320018621349" },
  @{ Cell = "D14"; Value = "320018621371" },
  @{ Cell = "D15"; Value = "320018621393" },
  @{ Cell = "D16"; Value = "320018621420" },
  @{ Cell = "D17"; Value = "320018621441" }
)

$scratch = $ws.Range("BA200")

foreach ($u in $updates) {
  $scratch.Formula = "=TEXT(" + $u.Value + ",""0"")"
  $scratch.Copy()
  $ws.Range($u.Cell).PasteSpecial(-4163)
}

$scratch.Clear()
$excel.CutCopyMode = $false
